# Atualização de bases das ligas, do dia: 10-06-2024 às 07:08
#
# Algumas partidas da "Ecuador LigaPro Serie A" estavam com os dados (id,
# times, placares, odds, etc.) trocados entre linhas que compartilhavam a
# mesma data. Este script corrige as colunas B e E..AD das linhas afetadas,
# mantendo A (id sequencial), C (Div) e D (Date) inalterados.
#
# Cada entrada do array abaixo representa uma linha da planilha, na ordem:
#   Row, B(id), E(HomeTeam), F(AwayTeam), G(FTHG), H(FTAG), I(HTHG), J(HTAG),
#   K(FTR), L(oddH_op), M(oddD_op), N(oddA_op), O(oddH), P(oddD), Q(oddA),
#   R(Ah), S(oddAHH), T(oddAHA), U(AhOU), V(oddAHOver), W(oddAHUnder),
#   X(PLH), Y(PLD), Z(PLA), AA(PL_Ahh), AB(PL_Aha), AC(PL_AhOver), AD(PL_AhUnder)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(139, 7528859, "Club Atletico Libertad", "Cumbaya FC", 3, 1, 2, 0, "H",
        1.727, 3.5, 4.333, 1.4, 4.2, 7,
        -1.25, 2, 1.8, 2.5, 1.95, 1.85,
        0.3999999999999999, -1, -1, 1, -1, 0.95, -1),

    @(140, 7528849, "Guayaquil City", "Gualaceo SC", 0, 2, 0, 1, "A",
        1.833, 3.5, 3.75, 2.15, 3.4, 3,
        -0.25, 1.825, 1.975, 2.5, 1.85, 1.95,
        -1, -1, 2, -1, 0.9750000000000001, -1, 0.95),

    @(142, 7528848, "Emelec", "Deportivo Cuenca", 2, 1, 0, 1, "H",
        1.75, 3.5, 4.2, 2.4, 3.1, 2.75,
        -0.25, 2.05, 1.75, 2.25, 1.8, 2,
        1.4, -1, -1, 1.05, -1, 0.8, -1),

    @(143, 7528858, "Orense", "SD Aucas", 1, 2, 1, 1, "A",
        2.2, 3.2, 3.2, 1.95, 3.2, 3.8,
        -0.5, 1.95, 1.85, 2.25, 1.85, 1.95,
        -1, -1, 2.8, -1, 0.8500000000000001, 0.8500000000000001, -1),

    @(144, 7528852, "Delfin SC", "Tecnico Universitario", 2, 2, 1, 0, "D",
        2.1, 3.4, 3.1, 2.1, 3.4, 3.1,
        -0.25, 1.8, 2, 2.25, 1.9, 1.9,
        -1, 2.4, -1, -0.5, 0.5, 0.8999999999999999, -1),

    @(145, 7528857, "Universidad Catolica del Ecuador", "Barcelona Guayaquil", 0, 1, 0, 0, "A",
        1.533, 4, 5.5, 1.5, 4.333, 5.25,
        -1, 1.8, 2, 3, 1.975, 1.825,
        -1, -1, 4.25, -1, 1, -1, 0.825)
)

foreach ($r in $rows) {
    $row = $r[0]

    $ws.Cells.Item($row, 2).Value2  = $r[1]   # B  id
    $ws.Cells.Item($row, 5).Value2  = $r[2]   # E  HomeTeam
    $ws.Cells.Item($row, 6).Value2  = $r[3]   # F  AwayTeam
    $ws.Cells.Item($row, 7).Value2  = $r[4]   # G  FTHG
    $ws.Cells.Item($row, 8).Value2  = $r[5]   # H  FTAG
    $ws.Cells.Item($row, 9).Value2  = $r[6]   # I  HTHG
    $ws.Cells.Item($row, 10).Value2 = $r[7]   # J  HTAG
    $ws.Cells.Item($row, 11).Value2 = $r[8]   # K  FTR
    $ws.Cells.Item($row, 12).Value2 = $r[9]   # L  oddH_op
    $ws.Cells.Item($row, 13).Value2 = $r[10]  # M  oddD_op
    $ws.Cells.Item($row, 14).Value2 = $r[11]  # N  oddA_op
    $ws.Cells.Item($row, 15).Value2 = $r[12]  # O  oddH
    $ws.Cells.Item($row, 16).Value2 = $r[13]  # P  oddD
    $ws.Cells.Item($row, 17).Value2 = $r[14]  # Q  oddA
    $ws.Cells.Item($row, 18).Value2 = $r[15]  # R  Ah
    $ws.Cells.Item($row, 19).Value2 = $r[16]  # S  oddAHH
    $ws.Cells.Item($row, 20).Value2 = $r[17]  # T  oddAHA
    $ws.Cells.Item($row, 21).Value2 = $r[18]  # U  AhOU
    $ws.Cells.Item($row, 22).Value2 = $r[19]  # V  oddAHOver
    $ws.Cells.Item($row, 23).Value2 = $r[20]  # W  oddAHUnder
    $ws.Cells.Item($row, 24).Value2 = $r[21]  # X  PLH
    $ws.Cells.Item($row, 25).Value2 = $r[22]  # Y  PLD
    $ws.Cells.Item($row, 26).Value2 = $r[23]  # Z  PLA
    $ws.Cells.Item($row, 27).Value2 = $r[24]  # AA PL_Ahh
    $ws.Cells.Item($row, 28).Value2 = $r[25]  # AB PL_Aha
    $ws.Cells.Item($row, 29).Value2 = $r[26]  # AC PL_AhOver
    $ws.Cells.Item($row, 30).Value2 = $r[27]  # AD PL_AhUnder
}
